$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = "'5.404"
$ws.Range("D5").Value = "'0.05991"
$ws.Range("D6").Value = "'3.391"
$ws.Range("D7").Value = "'6.400"
$ws.Range("D8").Value = "'0.8086"
$ws.Range("D9").Value = "'0.9571"
$ws.Range("D11").Value = "'0.07395"
$ws.Range("D12").Value = "'0.03396"
$ws.Range("D13").Value = "'0.03064"
$ws.Range("D14").Value = "'0.09415"
$ws.Range("D15").Value = "'3.997"
$ws.Range("D16").Value = "'0.001600"
$ws.Range("D17").Value = "'0.04798"
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D18").Value = "'0.006125"
$ws.Range("E18").Value = "17TigerCashTCH"
$ws.Range("B19").Value = "HotbitToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D19").Value = "'0.005067"
$ws.Range("E19").Value = "18HotbitTokenHTB"
$ws.Range("B20").Value = "BitKan"
$ws.Range("C20").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D20").Value = "'0.0009893"
$ws.Range("E20").Value = "19BitKanKAN"
$ws.Range("B21").Value = "NitroEx"
$ws.Range("C21").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D21").Value = "'0.0001000"
$ws.Range("E21").Value = "20NitroExNTX"
$ws.Range("B22").Value = "LEO"
$ws.Range("C22").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D22").Value = "'3.704"
$ws.Range("E22").Value = "21LEOLEO"
$ws.Range("B23").Value = "BTSEToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D23").Value = "'2.186"
$ws.Range("E23").Value = "22BTSETokenBTSE"
$ws.Range("B24").Value = "One"
$ws.Range("C24").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D24").Value = "'0.01116"
$ws.Range("E24").Value = "23OneONEBestin24h"
$ws.Range("D26").Value = "'0.1285"
$ws.Range("D40").Value = "'0.04014"
$ws.Range("D41").Value = "'0.006572"
$ws.Range("D43").Value = "'0.002901"
$ws.Range("D44").Value = "'0.005312"
$ws.Range("D45").Value = "'0.00005258"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"
$ws.Range("D48").Value = "'0.02469"
